$wb = $excel.ActiveWorkbook

$wsMS = $wb.Worksheets.Item("MSData")
$wsFlux = $wb.Worksheets.Item("FluxData")
$wsTracer = $wb.Worksheets.Item("TracerData")

# ----------------------------------------------------------------------
# MSData sheet: zoom 65 -> 100, selection A2:H34 -> A2, no longer the
# selected/active tab.
# ----------------------------------------------------------------------
$wsMS.Activate()
$excel.ActiveWindow.Zoom = 100
$wsMS.Range("A2").Select() | Out-Null

# ----------------------------------------------------------------------
# FluxData sheet: a new row is inserted as row 3 ("EX_glc__D_e.f"),
# pushing the previous rows 3-37 down to 4-38. A couple of values are
# also updated.
# ----------------------------------------------------------------------
$wsFlux.Activate()

$wsFlux.Rows.Item(3).Insert()

$wsFlux.Range("A3").Value = "EX_glc__D_e.f"
$wsFlux.Range("B3").Value = 3.73499001440729
$wsFlux.Range("C3").Value = 1.03731422547792
$wsFlux.Rows.Item(3).RowHeight = 13.8

# fix up the BIOMASS.f row's relative fit penalty value
$wsFlux.Range("C2").Value = 0.0001

# the row that used to be row 3 (EX_c5sugal_e.f) is now row 4; update its values
$wsFlux.Range("B4").Value = 0.030092592592593
$wsFlux.Range("C4").Value = 0.030855319651359

$excel.ActiveWindow.Zoom = 100
$wsFlux.Range("A9").Select() | Out-Null

# ----------------------------------------------------------------------
# TracerData sheet: zoom 65 -> 100, selection simplified to A1.
# ----------------------------------------------------------------------
$wsTracer.Activate()
$excel.ActiveWindow.Zoom = 100
$wsTracer.Range("A1").Select() | Out-Null

# ----------------------------------------------------------------------
# FluxData ends up as the active/selected sheet (activeTab = 1).
# ----------------------------------------------------------------------
$wsFlux.Activate()
